# Update countries & provincias Spain
# Applies the 19-May-2020 23:35 data refresh to the "Pais" sheet:
#  - refreshes totals for a handful of existing-position countries
#  - re-sorts the Benin..Mauritania block (Benin moves to the bottom,
#    the other 25 countries shift up one row and pick up fresh totals)
#  - bumps the "last updated" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 23:35"

# --- Countries whose row position is unchanged, only totals refreshed ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1566956
$ws.Range("C4").Value = 16662
$ws.Range("D4").Value = 362242
$ws.Range("E4").Value = 1111353
$ws.Range("G4").Value = 1380
$ws.Range("H4").Value = 93361

# Row 7: Brasil
$ws.Range("B7").Value = 265896
$ws.Range("C7").Value = 10528
$ws.Range("E7").Value = 147597
$ws.Range("G7").Value = 987
$ws.Range("H7").Value = 17840

# Row 28: Suiza
$ws.Range("D28").Value = 27700
$ws.Range("E28").Value = 1027

# Row 97: Eslovenia
$ws.Range("D97").Value = 1338
$ws.Range("E97").Value = 25

# --- Rows 138-163: Benin..Mauritania block re-sorted ---
# "Togo".."Mauritania" move up to occupy rows 138-162 (each keeping its own
# refreshed totals), and "Benin" drops to row 163 with its refreshed totals.

$countryBlock = @(
    @{Row=138; Country="Togo";                   B=338; C=8;   D=107; E=219; G=0; H=12},
    @{Row=139; Country="Cabo Verde";              B=335; C=7;   D=85;  E=247; G=0; H=3},
    @{Row=140; Country="Isla de Man";             B=335; C=0;   D=300; E=11;  G=0; H=24},
    @{Row=141; Country="Mauricio";                B=332; C=0;   D=322; E=0;   G=0; H=10},
    @{Row=142; Country="Madagascar";              B=326; C=4;   D=119; E=205; G=1; H=2},
    @{Row=143; Country="Vietnam";                 B=324; C=0;   D=263; E=61;  G=0; H=0},
    @{Row=144; Country="Montenegro";               B=324; C=0;   D=312; E=3;   G=0; H=9},
    @{Row=145; Country="Ruanda";                  B=308; C=11;  D=209; E=99;  G=0; H=0},
    @{Row=146; Country="Sudan del Sur";            B=290; C=0;   D=4;   E=282; G=0; H=4},
    @{Row=147; Country="Uganda";                  B=260; C=12;  D=63;  E=197; G=0; H=0},
    @{Row=148; Country="Nicaragua";                B=254; C=229; D=199; E=38;  G=9; H=17},
    @{Row=149; Country="Santo Tome y Principe";    B=246; C=0;   D=4;   E=235; G=0; H=7},
    @{Row=150; Country="Liberia";                 B=233; C=4;   D=125; E=85;  G=1; H=23},
    @{Row=151; Country="Guayana Francesa";         B=210; C=0;   D=131; E=78;  G=0; H=1},
    @{Row=152; Country="Suazilandia";              B=208; C=3;   D=87;  E=119; G=0; H=2},
    @{Row=153; Country="Birmania";                 B=193; C=5;   D=104; E=83;  G=0; H=6},
    @{Row=154; Country="Martinica";                B=192; C=0;   D=91;  E=87;  G=0; H=14},
    @{Row=155; Country="Islas Feroe";              B=187; C=0;   D=187; E=0;   G=0; H=0},
    @{Row=156; Country="Yemen";                   B=167; C=37;  D=5;   E=134; G=8; H=28},
    @{Row=157; Country="Guadalupe";                B=155; C=0;   D=109; E=33;  G=0; H=13},
    @{Row=158; Country="Gibraltar";                B=147; C=0;   D=145; E=2;   G=0; H=0},
    @{Row=159; Country="Mozambique";               B=146; C=1;   D=44;  E=102; G=0; H=0},
    @{Row=160; Country="Brunei";                  B=141; C=0;   D=136; E=4;   G=0; H=1},
    @{Row=161; Country="Mongolia";                 B=140; C=0;   D=26;  E=114; G=0; H=0},
    @{Row=162; Country="Mauritania";               B=131; C=50;  D=7;   E=120; G=0; H=4},
    @{Row=163; Country="Benin";                   B=130; C=0;   D=83;  E=45;  G=0; H=2}
)

foreach ($item in $countryBlock) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Country
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("G$r").Value = $item.G
    $ws.Range("H$r").Value = $item.H
}

Write-Host "Update complete"
